# case with 380 kV done
# Updates the loading_percent results table (rows 2-25, columns B,D,E,F,G,H,I,J,M,O)
# with the recomputed per-line loading percentages for the 380 kV case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B=10.25453685469846;  D=9.917401077243495;  E=13.40927253346448;  F=25.98016527116568;  G=24.73266510721236;  H=12.89130559803428;  I=24.25145179729221;  J=9.546370430325585;  M=25.8758476722981;   O=19.14433683646001 }
    3  = @{ B=9.787284860100545;  D=9.98024437576033;    E=13.52716516042191;  F=26.01585284262192;  G=24.51050607050577;  H=12.92163636755904;  I=23.49559474754147;  J=9.607178246877176;  M=24.72276231375537;  O=19.14645055845156 }
    4  = @{ B=9.487623971978266;  D=10.02081962214527;   E=13.60312101450734;  F=26.04717159191241;  G=24.38466145936372;  H=12.94325875319986;  I=23.02465761465173;  J=9.646448780127299;  M=23.98281939100087;  O=19.15389177293605 }
    5  = @{ B=9.362416834093104;  D=10.03785549095546;   E=13.63497352494129;  F=26.06228722265149;  G=24.33609202868887;  H=12.95282184824544;  I=22.83136587838137;  J=9.662939282163412;  M=23.67354116590206;  O=19.15846488050084 }
    6  = @{ B=9.34144327770055;   D=10.04071457813833;   E=13.64031703019994;  F=26.06493888478447;  G=24.32819252226616;  H=12.95445512727168;  I=22.79919688291642;  J=9.665706988823818;  M=23.62172722070081;  O=19.15931715705664 }
    7  = @{ B=9.485947739076757;  D=10.02104734356728;   E=13.60354694192532;  F=26.04736593698521;  G=24.38399537868633;  H=12.94338468372297;  I=23.02205594154593;  J=9.646669201803519;  M=23.97867930564853;  O=19.15394721532058 }
    8  = @{ B=10.09614276322378;  D=9.938657144880827;   E=13.44918252087398;  F=25.99051109722327;  G=24.65390915113088;  H=12.9011397911988;   I=23.99244694396022;  J=9.566936133821336;  M=25.48505250624957;  O=19.14378851558847 }
    9  = @{ B=11.18716782704496;  D=9.792830653800857;   E=13.17468674260087;  F=25.95413147386733;  G=25.26414830156611;  H=12.84219451920232;  I=25.82751908155621;  J=9.425879579637506;  M=28.1750978223021;   O=19.17273841404601 }
    10 = @{ B=11.91973579520242;  D=9.69522637017756;    E=12.99006049213558;  F=25.97373273387927;  G=25.75756107637644;  H=12.81358269592451;  I=27.11796810969694;  J=9.331505259525088;  M=29.97909531598403;  O=19.2238924956592 }
    11 = @{ B=12.23729015221987;  D=9.652881090273716;   E=12.90973824979282;  F=25.99278696897225;  G=25.99080033209485;  H=12.80378087275235;  I=27.68952754822376;  J=9.290568714630906;  M=30.76055110595963;  O=19.25364499707432 }
    12 = @{ B=12.3552359802385;   D=9.637140701295785;   E=12.87984700987434;  F=26.00146327006744;  G=26.08029485849865;  H=12.80053290812302;  I=27.90352491078817;  J=9.275352984595015;  M=31.05071543406407;  O=19.26584068286281 }
    13 = @{ B=12.32993743241465;  D=9.640517573321947;   E=12.8862613028028;   F=25.99952967222847;  G=26.06096996876369;  H=12.80121176535435;  I=27.85754843311305;  J=9.278617254501647;  M=30.98848104427516;  O=19.26317286720039 }
    14 = @{ B=12.24704006978876;  D=9.651580213790592;   E=12.9072685679988;   F=25.99347148463165;  G=25.9981401453579;   H=12.80350435748751;  I=27.70718328464474;  J=9.289311178946265;  M=30.78453907223427;  O=19.25462974627964 }
    15 = @{ B=12.19596156520551;  D=9.658394777785427;   E=12.92020444297663;  F=25.98995097725521;  G=25.95980492949553;  H=12.80496907794584;  I=27.61475657993344;  J=9.295898745531332;  M=30.65886591818147;  O=19.24951771410526 }
    16 = @{ B=11.89866272742029;  D=9.698035032122952;   E=12.99538332046924;  F=25.97269179921721;  G=25.74248735450391;  H=12.81428808270753;  I=27.08028510374492;  J=9.334220630554608;  M=29.92722592227451;  O=19.22207830723065 }
    17 = @{ B=11.71222255569353;  D=9.722878966002881;   E=13.04244043006747;  F=25.96470316118587;  G=25.61135570574974;  H=12.82082935641123;  I=26.74828088106421;  J=9.358240204960957;  M=29.46825790757741;  O=19.20690362484181 }
    18 = @{ B=11.60351212794432;  D=9.737362044807362;   E=13.069851534176;    F=25.96106243874606;  G=25.5367646324955;   H=12.82489414160925;  I=26.55587580389773;  J=9.372243435303393;  M=29.20058841996738;  O=19.19878596097469 }
    19 = @{ B=11.56645293568122;  D=9.742299023265108;   E=13.07919179573613;  F=25.95999349822867;  G=25.51165507198945;  H=12.82632229552897;  I=26.4904895024457;   J=9.377016969330356;  M=29.10933155943864;  O=19.19614236974913 }
    20 = @{ B=11.73222246210769;  D=9.720214265002539;   E=13.03739541847896;  F=25.96545478744876;  G=25.62522939697851;  H=12.82010171508478;  I=26.78377448402868;  J=9.355663848858194;  M=29.51749797626599;  O=19.20845583823905 }
    21 = @{ B=12.2714519360971;   D=9.648322851575415;   E=12.90108399097948;  F=25.99521125852484;  G=26.01656372176711;  H=12.8028183690301;   I=27.75141698525725;  J=9.286162356986875;  M=30.84459884260012;  O=19.25711388493515 }
    22 = @{ B=12.61041712741359;  D=9.603055949783359;   E=12.81505609028351;  F=26.02317396376178;  G=26.27911009383696;  H=12.79422652993699;  I=28.36952289460204;  J=9.242406089882566;  M=31.67834415120953;  O=19.29432760401874 }
    23 = @{ B=12.43074966401699;  D=9.627058753611992;   E=12.86069150447645;  F=26.00747009557275;  G=26.13839388128558;  H=12.79856427399851;  I=28.04100199104799;  J=9.265607360980081;  M=31.23646659794303;  O=19.27397203543869 }
    24 = @{ B=11.72318524601403;  D=9.721418353351877;   E=13.03967515486932;  F=25.9651120118191;   G=25.6189546076897;   H=12.82042973449486;  I=26.76773257959455;  J=9.356828014858792;  M=29.4952483763353;   O=19.20775219318957 }
    25 = @{ B=10.90385125002249;  D=9.830601873805893;   E=13.24594076963291;  F=25.95586908022578;  G=25.09081731789951;  H=12.85556860590771;  I=25.34014337090489;  J=9.462408222972616;  M=27.47693853303905;  O=19.15966005196335 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
